$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 83
$ws.Cells.Item(3, 6).Value = 66
$ws.Cells.Item(4, 6).Value = 960
$ws.Cells.Item(5, 6).Value = 1261
$ws.Cells.Item(6, 6).Value = 1741
$ws.Cells.Item(7, 6).Value = 917
$ws.Cells.Item(8, 6).Value = 572
$ws.Cells.Item(9, 6).Value = 2590
$ws.Cells.Item(10, 6).Value = 739
$ws.Cells.Item(11, 6).Value = 574
$ws.Cells.Item(12, 6).Value = 570
$ws.Cells.Item(13, 6).Value = 33
$ws.Cells.Item(14, 6).Value = 648
$ws.Cells.Item(15, 6).Value = 336
$ws.Cells.Item(16, 6).Value = 290
$ws.Cells.Item(18, 6).Value = 2126
$ws.Cells.Item(20, 6).Value = 708
$ws.Cells.Item(21, 6).Value = 6
$ws.Cells.Item(22, 6).Value = 2621
$ws.Cells.Item(23, 6).Value = 3
$ws.Cells.Item(28, 6).Value = 411
$ws.Cells.Item(29, 6).Value = 1776
$ws.Cells.Item(32, 6).Value = 526
$ws.Cells.Item(33, 6).Value = 547
$ws.Cells.Item(34, 6).Value = 210
$ws.Cells.Item(36, 6).Value = 338
$ws.Cells.Item(37, 6).Value = 4584
$ws.Cells.Item(38, 6).Value = 152

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(2, 6).Value = 387
$ws.Cells.Item(4, 6).Value = 4201
$ws.Cells.Item(11, 6).Value = 2
$ws.Cells.Item(12, 6).Value = 13
$ws.Cells.Item(14, 6).Value = 325
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(22, 6).Value = 60
$ws.Cells.Item(25, 6).Value = 1764
$ws.Cells.Item(28, 6).Value = 14
$ws.Cells.Item(29, 6).Value = 268
$ws.Cells.Item(39, 6).Value = 3

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(3, 6).Value = 63
$ws.Cells.Item(4, 6).Value = 1422
$ws.Cells.Item(6, 6).Value = 531
$ws.Cells.Item(7, 6).Value = 131
$ws.Cells.Item(8, 6).Value = 202

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 63
$ws.Cells.Item(3, 6).Value = 1422
$ws.Cells.Item(4, 6).Value = 531
$ws.Cells.Item(5, 6).Value = 83
$ws.Cells.Item(6, 6).Value = 66
$ws.Cells.Item(7, 6).Value = 960
$ws.Cells.Item(8, 6).Value = 1261
$ws.Cells.Item(9, 6).Value = 1741
$ws.Cells.Item(12, 6).Value = 131
$ws.Cells.Item(14, 6).Value = 917
$ws.Cells.Item(15, 6).Value = 572
$ws.Cells.Item(16, 6).Value = 2590
$ws.Cells.Item(17, 6).Value = 739
$ws.Cells.Item(18, 6).Value = 574
$ws.Cells.Item(19, 6).Value = 570
$ws.Cells.Item(20, 6).Value = 33
$ws.Cells.Item(21, 6).Value = 648
$ws.Cells.Item(22, 6).Value = 336
$ws.Cells.Item(24, 6).Value = 290
$ws.Cells.Item(25, 6).Value = 325
$ws.Cells.Item(27, 6).Value = 2126
$ws.Cells.Item(29, 6).Value = 708
$ws.Cells.Item(31, 6).Value = 6
$ws.Cells.Item(32, 6).Value = 2621
$ws.Cells.Item(36, 6).Value = 520
$ws.Cells.Item(38, 6).Value = 202
$ws.Cells.Item(40, 6).Value = 411
$ws.Cells.Item(41, 6).Value = 411
$ws.Cells.Item(43, 6).Value = 14
$ws.Cells.Item(44, 6).Value = 526
$ws.Cells.Item(45, 6).Value = 547
$ws.Cells.Item(47, 6).Value = 338
$ws.Cells.Item(48, 6).Value = 4584
$ws.Cells.Item(49, 6).Value = 152
